# Rebalance stats and elements
# Updates the HP/EP table counts+totals and repositions the saved view/selection
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    ,("C2", 35)
    ,("C3", 32)
    ,("C4", 30)
    ,("C5", 29)
    ,("C6", 28)
    ,("C7", 27)
    ,("C8", 26)
    ,("C9", 26)
    ,("C10", 25)
    ,("C11", 24)
    ,("C12", 23)
    ,("C13", 22)
    ,("C14", 22)
    ,("C15", 21)
    ,("C16", 20)
    ,("C19", 25)
    ,("D19", 601)
    ,("G19", 24)
    ,("H19", 574)
    ,("C20", 18)
    ,("G20", 16)
    ,("C21", 17)
    ,("G21", 15)
    ,("C22", 17)
    ,("G22", 15)
    ,("C23", 17)
    ,("G23", 15)
    ,("C24", 17)
    ,("G24", 15)
    ,("C25", 17)
    ,("G25", 14)
    ,("C26", 16)
    ,("G26", 14)
    ,("C27", 16)
    ,("G27", 14)
    ,("C28", 16)
    ,("G28", 14)
    ,("C29", 15)
    ,("G29", 13)
    ,("C30", 14)
    ,("G30", 12)
    ,("H30", 198)
    ,("C31", 13)
    ,("G31", 12)
    ,("H31", 187)
    ,("C32", 12)
    ,("G32", 12)
    ,("H32", 171)
    ,("C33", 12)
    ,("G33", 12)
    ,("H33", 165)
    ,("C36", 18)
    ,("G36", 19)
    ,("C37", 17)
    ,("G37", 18)
    ,("C38", 17)
    ,("G38", 18)
    ,("C39", 17)
    ,("G39", 18)
    ,("C40", 15)
    ,("G40", 17)
    ,("C41", 15)
    ,("G41", 16)
    ,("C42", 15)
    ,("G42", 15)
    ,("C43", 15)
    ,("G43", 14)
    ,("C44", 15)
    ,("G44", 14)
    ,("C45", 14)
    ,("G45", 14)
    ,("C46", 13)
    ,("G46", 13)
    ,("C47", 11)
    ,("G47", 13)
    ,("C48", 10)
    ,("G48", 13)
    ,("C49", 10)
    ,("G49", 13)
    ,("C50", 9)
    ,("G50", 12)
)

foreach ($pair in $updates) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# Restore the scroll position / active selection recorded in the saved view
$ws.Range("A16:M16").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G34").Select() | Out-Null
